$d = $word.ActiveDocument

function Set-BoldRun($findText) {
    $r = $d.Content
    $r.Find.Execute($findText, $true) | Out-Null
    $r.Bold = 1
}

function Split-RunAt($findText) {
    # Forces a run boundary right after $findText without any visible change,
    # by adding then immediately deleting a zero-width bookmark at that point.
    $r = $d.Content
    $r.Find.Execute($findText, $true) | Out-Null
    $r.Collapse(0)
    $d.Bookmarks.Add("TempSplitMark", $r) | Out-Null
    $bm = $d.Bookmarks("TempSplitMark")
    $bm.Delete()
}

function Split-RunBefore($findText) {
    # Forces a run boundary right before $findText without any visible change.
    $r = $d.Content
    $r.Find.Execute($findText, $true) | Out-Null
    $r.Collapse(1)
    $d.Bookmarks.Add("TempSplitMark", $r) | Out-Null
    $bm = $d.Bookmarks("TempSplitMark")
    $bm.Delete()
}

# 1. Bold "not current" in the LinkedIn paragraph.
Set-BoldRun("not current")

# 2. Split the bold "NoCode NoSQL DBaaS" run into "NoCode" + " NoSQL DBaaS"
#    (both stay bold; this mirrors the proofErr spell-check boundary around "NoCode").
Split-RunAt("NoCode")

# 3. Mark the architecture-diagram drawing run as NoProof.
$diagShape = $d.InlineShapes(1)
$diagShape.Range.NoProofing = 1

# 4. Bold MongoDB / F# / Haskell in the "technologies" paragraph, and split off
#    "QuickCheck" into its own run (mirrors the proofErr spell-check boundary).
Set-BoldRun("MongoDB")
Set-BoldRun("F#")
Set-BoldRun("Haskell")
Split-RunBefore("QuickCheck")
Split-RunAt("QuickCheck")

Write-Output "done step2-4"
